$d = $word.ActiveDocument

# --- Paragraph 1 (heading): merge the two split runs "Синтетизир" + "ани звуци"
# into a single run reading "Синтетизирани звуци", and drop the _GoBack bookmark
# that used to sit between them (it is re-added at the end of paragraph 3 below).
$heading = $d.Paragraphs.Item(26)
$headingXml = @'
<w:p w14:paraId="108A740C" w14:textId="47B3A619" w:rsidR="00811966" w:rsidRDefault="00811966" w:rsidP="009A7F8C" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t>Синтетизирани звуци</w:t></w:r></w:p>
'@
$heading.Range.InsertXML($headingXml)

# --- Paragraph 3 (body): split "clap.wav"/"kick.wav" into clap/Synth/.wav runs
# and append the large new block of text describing how the sounds were made,
# ending with the _GoBack bookmark (its new location after the edit).
$body = $d.Paragraphs.Item(28)
$bodyXml = @'
<w:p w14:paraId="2ED5B83C" w14:textId="690A983F" w:rsidR="00811966" w:rsidRPr="00811966" w:rsidRDefault="00811966" w:rsidP="009A7F8C" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">Синтетизирани беа звуците </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>clap</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>Synth</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">.wav </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">и </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>kick</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>Synth</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>.wav</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">. Целта беше да звучат што повеќе на звуците снимени во </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>FL Studio.</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">За снимање на истите е искористена </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">Play </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">функција која е мапирана на копче по број 121. Конкретно </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>clapSynth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">е направено со користење на два </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>WhiteNoise</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">-а кои се обликувани со две енвелопи така што првиата енвелопа и </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>WhiteNoise</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">се користат за првиот дел од звукот и вторите соодветно за вториот. Така што кога ќе се притисне копчето мапирано на 121, ќе почне да се повтара звукот бесконечен број на пати, се додека не се притисне повторно, со што се запира звукот. Додека </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>kickSynth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">е направен со </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t>SinOsc</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve"> и повеќе енвелопи со кои се одредува фреквенцијата и амплитудата на </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t xml:space="preserve">звуците кои се спојуваат да го направат </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/></w:rPr><w:t>kickSynth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="36"/><w:lang w:val="mk-MK"/></w:rPr><w:t>. За самата употреба на звуците повторно мапирањето и пуштањето на звукот е исто.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$body.Range.InsertXML($bodyXml)
